$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "249.15"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.54"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.410"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05712"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.412"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8135"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9257"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1423"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07430"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03133"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03055"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09355"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.756"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001574"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04771"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.01829"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005796"

$ws.Range("E19").Value = "18OneONEWorstin24h"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006479"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.005005"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.001025"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.700"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.159"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03993"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006874"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1065"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007515"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005897"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5006"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
